# Updated symbol list (crypto prices) - Price column (D) refresh.
# Source values are stored as text (coinranking.com scrape), so each new
# value is written with a leading apostrophe to force Excel to keep the
# cell as text instead of auto-converting the numeric-looking string to
# a number (which would also introduce binary floating point noise).

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value  = "'263.35"
$ws.Range("D3").Value  = "'21.64"
$ws.Range("D5").Value  = "'0.06179"
$ws.Range("D7").Value  = "'6.546"
$ws.Range("D8").Value  = "'1.392"
$ws.Range("D9").Value  = "'0.8248"
$ws.Range("D10").Value = "'0.1623"
$ws.Range("D11").Value = "'0.08202"
$ws.Range("D12").Value = "'0.03557"
$ws.Range("D13").Value = "'0.03181"
$ws.Range("D14").Value = "'0.09207"
$ws.Range("D15").Value = "'3.771"
$ws.Range("D16").Value = "'0.001626"
$ws.Range("D18").Value = "'0.006410"
$ws.Range("D21").Value = "'0.0001504"
$ws.Range("D22").Value = "'3.726"
$ws.Range("D23").Value = "'2.236"
$ws.Range("D24").Value = "'0.01358"
$ws.Range("D40").Value = "'0.04701"
$ws.Range("D41").Value = "'0.007000"
$ws.Range("D42").Value = "'0.1124"
$ws.Range("D43").Value = "'0.003548"
$ws.Range("D45").Value = "'0.00006087"
$ws.Range("D46").Value = "'0.0009906"
$ws.Range("D48").Value = "'0.9805"
